# Weekly data update: extend "Weekly Expenditure" sheet with Wk43 entries
# and fix Material Description (column D) alignment for prior rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix column D (Material Description) alignment for existing rows 666:771 ---
# (left + vertical-center, matching the style already used elsewhere in the column)
$ws.Range("D666:D771").HorizontalAlignment = -4131

# --- Prime formatting for the new rows by copying from an already-styled row ---
$ws.Range("D2").Copy()
$ws.Range("D772:D814").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("G772:G814").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Append the new Wk43 expenditure rows (772-814) ---
# Row 772
$ws.Cells.Item(772, 1).Value = 'Wk43'
$ws.Cells.Item(772, 2).Value = 'SAP'
$ws.Cells.Item(772, 3).Value = 'XS-PTS-0871'
$ws.Cells.Item(772, 4).Value = 'AFX-8320M-01-REV1 PPMU & FPGA I2C REV2'
$ws.Cells.Item(772, 5).Value = '1pcs'
$ws.Cells.Item(772, 6).Value = 'Fishes'
$ws.Cells.Item(772, 7).Value = '20/10/2025'
$ws.Cells.Item(772, 8).Value = 2242.2600000000002

# Row 773
$ws.Cells.Item(773, 1).Value = 'Wk43'
$ws.Cells.Item(773, 2).Value = 'SAP'
$ws.Cells.Item(773, 3).Value = 'XS-PTS-0976'
$ws.Cells.Item(773, 4).Value = '10618.334 MORAY FPGA INTERCONNECTION'
$ws.Cells.Item(773, 5).Value = '1pcs'
$ws.Cells.Item(773, 6).Value = 'Fishes'
$ws.Cells.Item(773, 7).Value = '20/10/2025'
$ws.Cells.Item(773, 8).Value = 237

# Row 774
$ws.Cells.Item(774, 1).Value = 'Wk43'
$ws.Cells.Item(774, 2).Value = 'SAP'
$ws.Cells.Item(774, 3).Value = 'XS-PTS-0416'
$ws.Cells.Item(774, 4).Value = '10105.095 SIMATIC S7-1200, DIGITAL I/O'
$ws.Cells.Item(774, 5).Value = '1pcs'
$ws.Cells.Item(774, 6).Value = 'Fishes'
$ws.Cells.Item(774, 7).Value = '20/10/2025'
$ws.Cells.Item(774, 8).Value = 246.94

# Row 775
$ws.Cells.Item(775, 1).Value = 'Wk43'
$ws.Cells.Item(775, 2).Value = 'SAP'
$ws.Cells.Item(775, 3).Value = 11156537
$ws.Cells.Item(775, 4).Value = 'PTS-1150 300-01698-010 Semiconductor'
$ws.Cells.Item(775, 5).Value = '600pcs'
$ws.Cells.Item(775, 6).Value = 'Fishes'
$ws.Cells.Item(775, 7).Value = '20/10/2025'
$ws.Cells.Item(775, 8).Value = 1404

# Row 776
$ws.Cells.Item(776, 1).Value = 'Wk43'
$ws.Cells.Item(776, 2).Value = 'SAP'
$ws.Cells.Item(776, 3).Value = 'XS-PTS-0876'
$ws.Cells.Item(776, 4).Value = 'HX 2067 HPN CRC CO CONTACT CLEANER'
$ws.Cells.Item(776, 5).Value = '5pcs'
$ws.Cells.Item(776, 6).Value = 'Fishes'
$ws.Cells.Item(776, 7).Value = '20/10/2025'
$ws.Cells.Item(776, 8).Value = 167.05

# Row 777
$ws.Cells.Item(777, 1).Value = 'Wk43'
$ws.Cells.Item(777, 2).Value = 'SAP'
$ws.Cells.Item(777, 3).Value = 11155949
$ws.Cells.Item(777, 4).Value = 'PTS-1146 70902.631 X2544 ALIGNER & NEST'
$ws.Cells.Item(777, 5).Value = '6pcs'
$ws.Cells.Item(777, 6).Value = 'Fishes'
$ws.Cells.Item(777, 7).Value = '20/10/2025'
$ws.Cells.Item(777, 8).Value = 17336.88

# Row 778
$ws.Cells.Item(778, 1).Value = 'Wk43'
$ws.Cells.Item(778, 2).Value = 'SAP'
$ws.Cells.Item(778, 3).Value = 11151246
$ws.Cells.Item(778, 4).Value = 'PTS-1069 TW.50.15.FI.0S.150.00 X2637 (Station A DOME ADAPTOR/CAP  (BOTTOM))'
$ws.Cells.Item(778, 5).Value = '5pcs'
$ws.Cells.Item(778, 6).Value = 'Fishes'
$ws.Cells.Item(778, 7).Value = '20/10/2025'
$ws.Cells.Item(778, 8).Value = 2777.8

# Row 779
$ws.Cells.Item(779, 1).Value = 'Wk43'
$ws.Cells.Item(779, 2).Value = 'SAP'
$ws.Cells.Item(779, 3).Value = 11155143
$ws.Cells.Item(779, 4).Value = 'PTS-1136 300-01519-015 Semiconductor'
$ws.Cells.Item(779, 5).Value = '1200pcs'
$ws.Cells.Item(779, 6).Value = 'Fishes'
$ws.Cells.Item(779, 7).Value = '20/10/2025'
$ws.Cells.Item(779, 8).Value = 2232

# Row 780
$ws.Cells.Item(780, 1).Value = 'Wk43'
$ws.Cells.Item(780, 2).Value = 'SAP'
$ws.Cells.Item(780, 3).Value = 11155949
$ws.Cells.Item(780, 4).Value = 'PTS-1146 70902.631 X2544 ALIGNER & NEST'
$ws.Cells.Item(780, 5).Value = '4pcs'
$ws.Cells.Item(780, 6).Value = 'Fishes'
$ws.Cells.Item(780, 7).Value = '20/10/2025'
$ws.Cells.Item(780, 8).Value = 11557.92

# Row 781
$ws.Cells.Item(781, 1).Value = 'Wk43'
$ws.Cells.Item(781, 2).Value = 'SAP'
$ws.Cells.Item(781, 3).Value = 11156537
$ws.Cells.Item(781, 4).Value = 'PTS-1150 300-01698-010 Semiconductor'
$ws.Cells.Item(781, 5).Value = '200pcs'
$ws.Cells.Item(781, 6).Value = 'Fishes'
$ws.Cells.Item(781, 7).Value = '20/10/2025'
$ws.Cells.Item(781, 8).Value = 468

# Row 782
$ws.Cells.Item(782, 1).Value = 'Wk43'
$ws.Cells.Item(782, 2).Value = 'SAP'
$ws.Cells.Item(782, 3).Value = 'XS-PTS-0108'
$ws.Cells.Item(782, 4).Value = 'CDUJB10-6D AIR CYLINDER'
$ws.Cells.Item(782, 5).Value = '5pcs'
$ws.Cells.Item(782, 6).Value = 'Fishes'
$ws.Cells.Item(782, 7).Value = '21/10/2025'
$ws.Cells.Item(782, 8).Value = 134.96

# Row 783
$ws.Cells.Item(783, 1).Value = 'Wk43'
$ws.Cells.Item(783, 2).Value = 'SAP'
$ws.Cells.Item(783, 3).Value = 11156541
$ws.Cells.Item(783, 4).Value = 'PTS-1154 70902.708 TOP NEST X2544'
$ws.Cells.Item(783, 5).Value = '4pcs'
$ws.Cells.Item(783, 6).Value = 'Fishes'
$ws.Cells.Item(783, 7).Value = '21/10/2025'
$ws.Cells.Item(783, 8).Value = 1568.64

# Row 784
$ws.Cells.Item(784, 1).Value = 'Wk43'
$ws.Cells.Item(784, 2).Value = 'SAP'
$ws.Cells.Item(784, 3).Value = 'XS-PTS-0416'
$ws.Cells.Item(784, 4).Value = '10105.095 SIMATIC S7-1200, DIGITAL I/O'
$ws.Cells.Item(784, 5).Value = '1pcs'
$ws.Cells.Item(784, 6).Value = 'Fishes'
$ws.Cells.Item(784, 7).Value = '21/10/2025'
$ws.Cells.Item(784, 8).Value = 246.94

# Row 785
$ws.Cells.Item(785, 1).Value = 'Wk43'
$ws.Cells.Item(785, 2).Value = 'SAP'
$ws.Cells.Item(785, 3).Value = 11156537
$ws.Cells.Item(785, 4).Value = 'PTS-1150 300-01698-010 Semiconductor'
$ws.Cells.Item(785, 5).Value = '200pcs'
$ws.Cells.Item(785, 6).Value = 'Fishes'
$ws.Cells.Item(785, 7).Value = '21/10/2025'
$ws.Cells.Item(785, 8).Value = 468

# Row 786
$ws.Cells.Item(786, 1).Value = 'Wk43'
$ws.Cells.Item(786, 2).Value = 'SAP'
$ws.Cells.Item(786, 3).Value = 11151246
$ws.Cells.Item(786, 4).Value = 'HX 2067 HPN CRC CO CONTACT CLEANER'
$ws.Cells.Item(786, 5).Value = '2pcs'
$ws.Cells.Item(786, 6).Value = 'Fishes'
$ws.Cells.Item(786, 7).Value = '21/10/2025'
$ws.Cells.Item(786, 8).Value = 1111.1199999999999

# Row 787
$ws.Cells.Item(787, 1).Value = 'Wk43'
$ws.Cells.Item(787, 2).Value = 'SAP'
$ws.Cells.Item(787, 3).Value = 11155143
$ws.Cells.Item(787, 4).Value = 'PTS-1136 300-01519-015 Semiconductor'
$ws.Cells.Item(787, 5).Value = '400pcs'
$ws.Cells.Item(787, 6).Value = 'Fishes'
$ws.Cells.Item(787, 7).Value = '21/10/2025'
$ws.Cells.Item(787, 8).Value = 744

# Row 788
$ws.Cells.Item(788, 1).Value = 'Wk43'
$ws.Cells.Item(788, 2).Value = 'SAP'
$ws.Cells.Item(788, 3).Value = 11156536
$ws.Cells.Item(788, 4).Value = 'PTS-1149 40Ways Ribbon Flat Cable (1.2M)'
$ws.Cells.Item(788, 5).Value = '10pcs'
$ws.Cells.Item(788, 6).Value = 'Fishes'
$ws.Cells.Item(788, 7).Value = '22/10/2025'
$ws.Cells.Item(788, 8).Value = 194.66

# Row 789
$ws.Cells.Item(789, 1).Value = 'Wk43'
$ws.Cells.Item(789, 2).Value = 'SAP'
$ws.Cells.Item(789, 3).Value = 11155143
$ws.Cells.Item(789, 4).Value = 'PTS-1136 300-01519-015 Semiconductor'
$ws.Cells.Item(789, 5).Value = '400pcs'
$ws.Cells.Item(789, 6).Value = 'Fishes'
$ws.Cells.Item(789, 7).Value = '22/10/2025'
$ws.Cells.Item(789, 8).Value = 744

# Row 790
$ws.Cells.Item(790, 1).Value = 'Wk43'
$ws.Cells.Item(790, 2).Value = 'SAP'
$ws.Cells.Item(790, 3).Value = 11151249
$ws.Cells.Item(790, 4).Value = 'PTS -1072 70192.696 LL RUBBER TIP X2637'
$ws.Cells.Item(790, 5).Value = '60pcs'
$ws.Cells.Item(790, 6).Value = 'Fishes'
$ws.Cells.Item(790, 7).Value = '22/10/2025'
$ws.Cells.Item(790, 8).Value = 2823.48

# Row 791
$ws.Cells.Item(791, 1).Value = 'Wk43'
$ws.Cells.Item(791, 2).Value = 'SAP'
$ws.Cells.Item(791, 3).Value = 11151237
$ws.Cells.Item(791, 4).Value = 'PTS-1060 70192.692 PNP RUBBER TIP'
$ws.Cells.Item(791, 5).Value = '600pcs'
$ws.Cells.Item(791, 6).Value = 'Fishes'
$ws.Cells.Item(791, 7).Value = '22/10/2025'
$ws.Cells.Item(791, 8).Value = 4288.43

# Row 792
$ws.Cells.Item(792, 1).Value = 'Wk43'
$ws.Cells.Item(792, 2).Value = 'SAP'
$ws.Cells.Item(792, 3).Value = 'XS-PTS-0108'
$ws.Cells.Item(792, 4).Value = 'CDUJB10-6D AIR CYLINDER'
$ws.Cells.Item(792, 5).Value = '1pcs'
$ws.Cells.Item(792, 6).Value = 'Fishes'
$ws.Cells.Item(792, 7).Value = '22/10/2025'
$ws.Cells.Item(792, 8).Value = 27

# Row 793
$ws.Cells.Item(793, 1).Value = 'Wk43'
$ws.Cells.Item(793, 2).Value = 'SAP'
$ws.Cells.Item(793, 3).Value = 11156541
$ws.Cells.Item(793, 4).Value = 'PTS-1154 70902.708 TOP NEST X2544'
$ws.Cells.Item(793, 5).Value = '4pcs'
$ws.Cells.Item(793, 6).Value = 'Fishes'
$ws.Cells.Item(793, 7).Value = '23/10/2025'
$ws.Cells.Item(793, 8).Value = 1568.64

# Row 794
$ws.Cells.Item(794, 1).Value = 'Wk43'
$ws.Cells.Item(794, 2).Value = 'SAP'
$ws.Cells.Item(794, 3).Value = 'XS-PTS-1028'
$ws.Cells.Item(794, 4).Value = '10618.397 X1767 SMU DAUGHTER CARD'
$ws.Cells.Item(794, 5).Value = '2pcs'
$ws.Cells.Item(794, 6).Value = 'Fishes'
$ws.Cells.Item(794, 7).Value = '23/10/2025'
$ws.Cells.Item(794, 8).Value = 2221.04

# Row 795
$ws.Cells.Item(795, 1).Value = 'Wk43'
$ws.Cells.Item(795, 2).Value = 'SAP'
$ws.Cells.Item(795, 3).Value = 11151242
$ws.Cells.Item(795, 4).Value = 'PTS-1065 TW.50.1A.FI.0S.136.00 X2637_SECONDARY BOARD'
$ws.Cells.Item(795, 5).Value = '2pcs'
$ws.Cells.Item(795, 6).Value = 'Fishes'
$ws.Cells.Item(795, 7).Value = '23/10/2025'
$ws.Cells.Item(795, 8).Value = 213.5

# Row 796
$ws.Cells.Item(796, 1).Value = 'Wk43'
$ws.Cells.Item(796, 2).Value = 'SAP'
$ws.Cells.Item(796, 3).Value = 11156542
$ws.Cells.Item(796, 4).Value = 'PTS-1155 TW.50.1A.FI.0S.123.02 SA30 -'
$ws.Cells.Item(796, 5).Value = '15pcs'
$ws.Cells.Item(796, 6).Value = 'Fishes'
$ws.Cells.Item(796, 7).Value = '23/10/2025'
$ws.Cells.Item(796, 8).Value = 2167.1999999999998

# Row 797
$ws.Cells.Item(797, 1).Value = 'Wk43'
$ws.Cells.Item(797, 2).Value = 'SAP'
$ws.Cells.Item(797, 3).Value = 11156543
$ws.Cells.Item(797, 4).Value = 'PTS-1156 TW.50.1A.FI.0S.089.01 SA30'
$ws.Cells.Item(797, 5).Value = '20pcs'
$ws.Cells.Item(797, 6).Value = 'Fishes'
$ws.Cells.Item(797, 7).Value = '23/10/2025'
$ws.Cells.Item(797, 8).Value = 1341.6

# Row 798
$ws.Cells.Item(798, 1).Value = 'Wk43'
$ws.Cells.Item(798, 2).Value = 'Expense'
$ws.Cells.Item(798, 3).Value = 'Expense'
$ws.Cells.Item(798, 4).Value = 'Wera Phillips Precision Screwdriver, PH0 Tip, 60 mm Blade, 157 mm Overall (ESD Safe Type)'
$ws.Cells.Item(798, 5).Value = '10pcs'
$ws.Cells.Item(798, 6).Value = 'Fishes'
$ws.Cells.Item(798, 7).Value = '24/10/2025'
$ws.Cells.Item(798, 8).Value = 138.72999999999999

# Row 799
$ws.Cells.Item(799, 1).Value = 'Wk43'
$ws.Cells.Item(799, 2).Value = 'Expense'
$ws.Cells.Item(799, 3).Value = 'Expense'
$ws.Cells.Item(799, 4).Value = 'Socket Precision Harden Brushing Guide'
$ws.Cells.Item(799, 5).Value = '10pcs'
$ws.Cells.Item(799, 6).Value = 'Fishes'
$ws.Cells.Item(799, 7).Value = '24/10/2025'
$ws.Cells.Item(799, 8).Value = 181.11

# Row 800
$ws.Cells.Item(800, 1).Value = 'Wk43'
$ws.Cells.Item(800, 2).Value = 'Expense'
$ws.Cells.Item(800, 3).Value = 'Expense'
$ws.Cells.Item(800, 4).Value = 'Beware Of Laser Sticker (10pcs Per Pack)'
$ws.Cells.Item(800, 5).Value = '2pack'
$ws.Cells.Item(800, 6).Value = 'Fishes'
$ws.Cells.Item(800, 7).Value = '24/10/2025'
$ws.Cells.Item(800, 8).Value = 30.06

# Row 801
$ws.Cells.Item(801, 1).Value = 'Wk43'
$ws.Cells.Item(801, 2).Value = 'Expense'
$ws.Cells.Item(801, 3).Value = 'Expense'
$ws.Cells.Item(801, 4).Value = 'Tolsen LED Telescopic Inspection Mirror 66006'
$ws.Cells.Item(801, 5).Value = '2pcs'
$ws.Cells.Item(801, 6).Value = 'Fishes'
$ws.Cells.Item(801, 7).Value = '24/10/2025'
$ws.Cells.Item(801, 8).Value = 20.04

# Row 802
$ws.Cells.Item(802, 1).Value = 'Wk43'
$ws.Cells.Item(802, 2).Value = 'SAP'
$ws.Cells.Item(802, 3).Value = 'XS-SPM-0081'
$ws.Cells.Item(802, 4).Value = 'HP-1810LS160-01 TIB PIN'
$ws.Cells.Item(802, 5).Value = '200pcs'
$ws.Cells.Item(802, 6).Value = 'Lisa'
$ws.Cells.Item(802, 7).Value = '24/10/2025'
$ws.Cells.Item(802, 8).Value = 200

# Row 803
$ws.Cells.Item(803, 1).Value = 'Wk43'
$ws.Cells.Item(803, 2).Value = 'SAP'
$ws.Cells.Item(803, 3).Value = 'XS-PTS-0867'
$ws.Cells.Item(803, 4).Value = '6K-57084-H062 Hyperspace Semicon Socket'
$ws.Cells.Item(803, 5).Value = '1pcs'
$ws.Cells.Item(803, 6).Value = 'Sihl'
$ws.Cells.Item(803, 7).Value = '24/10/2025'
$ws.Cells.Item(803, 8).Value = 190.99

# Row 804
$ws.Cells.Item(804, 1).Value = 'Wk43'
$ws.Cells.Item(804, 2).Value = 'SAP'
$ws.Cells.Item(804, 3).Value = 'XS-PTS-0502'
$ws.Cells.Item(804, 4).Value = 'Hyperspace Socket Model : 6K-76235-H04'
$ws.Cells.Item(804, 5).Value = '1pcs'
$ws.Cells.Item(804, 6).Value = 'Sihl'
$ws.Cells.Item(804, 7).Value = '24/10/2025'
$ws.Cells.Item(804, 8).Value = 191.35

# Row 805
$ws.Cells.Item(805, 1).Value = 'Wk43'
$ws.Cells.Item(805, 2).Value = 'SAP'
$ws.Cells.Item(805, 3).Value = 11156536
$ws.Cells.Item(805, 4).Value = 'PTS-1149 40Ways Ribbon Flat Cable (1.2M)'
$ws.Cells.Item(805, 5).Value = '8pcs'
$ws.Cells.Item(805, 6).Value = 'Fishes'
$ws.Cells.Item(805, 7).Value = '24/10/2025'
$ws.Cells.Item(805, 8).Value = 155.72999999999999

# Row 806
$ws.Cells.Item(806, 1).Value = 'Wk43'
$ws.Cells.Item(806, 2).Value = 'SAP'
$ws.Cells.Item(806, 3).Value = 11156537
$ws.Cells.Item(806, 4).Value = 'PTS-1150 300-01698-010 Semiconductor'
$ws.Cells.Item(806, 5).Value = '500pcs'
$ws.Cells.Item(806, 6).Value = 'Fishes'
$ws.Cells.Item(806, 7).Value = '24/10/2025'
$ws.Cells.Item(806, 8).Value = 1170

# Row 807
$ws.Cells.Item(807, 1).Value = 'Wk43'
$ws.Cells.Item(807, 2).Value = 'SAP'
$ws.Cells.Item(807, 3).Value = 'XS-PTS-1043'
$ws.Cells.Item(807, 4).Value = '10416.113 SAMTEC CABLE 1.2M X1767'
$ws.Cells.Item(807, 5).Value = '2pcs'
$ws.Cells.Item(807, 6).Value = 'Fishes'
$ws.Cells.Item(807, 7).Value = '24/10/2025'
$ws.Cells.Item(807, 8).Value = 1399.17

# Row 808
$ws.Cells.Item(808, 1).Value = 'Wk43'
$ws.Cells.Item(808, 2).Value = 'SAP'
$ws.Cells.Item(808, 3).Value = 'XS-PTS-0266'
$ws.Cells.Item(808, 4).Value = '70192.261 VACUUM PAD (For Empty Tray)'
$ws.Cells.Item(808, 5).Value = '40pcs'
$ws.Cells.Item(808, 6).Value = 'Fishes'
$ws.Cells.Item(808, 7).Value = '24/10/2025'
$ws.Cells.Item(808, 8).Value = 756.02

# Row 809
$ws.Cells.Item(809, 1).Value = 'Wk43'
$ws.Cells.Item(809, 2).Value = 'SAP'
$ws.Cells.Item(809, 3).Value = 11155143
$ws.Cells.Item(809, 4).Value = 'PTS-1136 300-01519-015 Semiconductor'
$ws.Cells.Item(809, 5).Value = '900pcs'
$ws.Cells.Item(809, 6).Value = 'Fishes'
$ws.Cells.Item(809, 7).Value = '24/10/2025'
$ws.Cells.Item(809, 8).Value = 930

# Row 810
$ws.Cells.Item(810, 1).Value = 'Wk43'
$ws.Cells.Item(810, 2).Value = 'SAP'
$ws.Cells.Item(810, 3).Value = 11151246
$ws.Cells.Item(810, 4).Value = 'PTS-1069 TW.50.15.FI.0S.150.00 X2637 (Station A DOME ADAPTOR/CAP  (BOTTOM))'
$ws.Cells.Item(810, 5).Value = '4pcs'
$ws.Cells.Item(810, 6).Value = 'Fishes'
$ws.Cells.Item(810, 7).Value = '24/10/2025'
$ws.Cells.Item(810, 8).Value = 2222.2399999999998

# Row 811
$ws.Cells.Item(811, 1).Value = 'Wk43'
$ws.Cells.Item(811, 2).Value = 'SAP'
$ws.Cells.Item(811, 3).Value = 11156537
$ws.Cells.Item(811, 4).Value = 'PTS-1150 300-01698-010 Semiconductor'
$ws.Cells.Item(811, 5).Value = '300pcs'
$ws.Cells.Item(811, 6).Value = 'Fishes'
$ws.Cells.Item(811, 7).Value = '24/10/2025'
$ws.Cells.Item(811, 8).Value = 702

# Row 812
$ws.Cells.Item(812, 1).Value = 'Wk43'
$ws.Cells.Item(812, 2).Value = 'SAP'
$ws.Cells.Item(812, 3).Value = 'XS-PTS-0933'
$ws.Cells.Item(812, 4).Value = 'X1629 70900.137 SEMICONDUCTOR PROBE PIN'
$ws.Cells.Item(812, 5).Value = '400pcs'
$ws.Cells.Item(812, 6).Value = 'Fishes'
$ws.Cells.Item(812, 7).Value = '24/10/2025'
$ws.Cells.Item(812, 8).Value = 4964

# Row 813
$ws.Cells.Item(813, 1).Value = 'Wk43'
$ws.Cells.Item(813, 2).Value = 'SAP'
$ws.Cells.Item(813, 3).Value = 'XS-PTS-0837'
$ws.Cells.Item(813, 4).Value = 'HX 0786 HPN VGA CABLE MALE TO MALE 5M'
$ws.Cells.Item(813, 5).Value = '1pcs'
$ws.Cells.Item(813, 6).Value = 'Fishes'
$ws.Cells.Item(813, 7).Value = '24/10/2025'
$ws.Cells.Item(813, 8).Value = 11.8

# Row 814
$ws.Cells.Item(814, 1).Value = 'Wk43'
$ws.Cells.Item(814, 2).Value = 'SAP'
$ws.Cells.Item(814, 3).Value = 11155947
$ws.Cells.Item(814, 4).Value = 'PTS-1144 800.403.00 X2544 MYCROFTL MP'
$ws.Cells.Item(814, 5).Value = '4pcs'
$ws.Cells.Item(814, 6).Value = 'Fishes'
$ws.Cells.Item(814, 7).Value = '24/10/2025'
$ws.Cells.Item(814, 8).Value = 1504.64

# --- Expand the AutoFilter range to cover the newly added rows ---
$ws.AutoFilterMode = $false
$ws.Range("A1:H814").AutoFilter()

# --- Keep the hidden _FilterDatabase defined name in sync with the filter range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Weekly Expenditure!_FilterDatabase") {
        $n.RefersTo = '=''Weekly Expenditure''!$A$1:$H$814'
    }
}

# --- Restore the active selection to reflect the new bottom of the sheet ---
$ws.Range("D823").Select()